$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing the old row 3 down to row 4
$ws.Rows.Item(3).Insert()

# --- Row 2 (existing row, values updated) ---
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Tff3"
$ws.Cells.Item(2,3).Value2 = "Ackr3"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 0.7901163333333333
$ws.Cells.Item(2,8).Value2 = 2.370349
$ws.Cells.Item(2,9).Value2 = 1
$ws.Cells.Item(2,10).Value2 = 1
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 9.688363000000001
$ws.Cells.Item(2,14).Value2 = 29.065089
$ws.Cells.Item(2,15).Value2 = 0.1053077753334822
$ws.Cells.Item(2,16).Value2 = 0.1053077753334822
$ws.Cells.Item(2,17).Value2 = 7.654933849562334
$ws.Cells.Item(2,18).Value2 = 68.89440464606101
$ws.Cells.Item(2,19).Value2 = 0.1053077753334822
$ws.Cells.Item(2,20).Value2 = 0.1053077753334822

# --- Row 3 (newly inserted row) ---
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Tff3"
$ws.Cells.Item(3,3).Value2 = "Ackr3"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 0.7901163333333333
$ws.Cells.Item(3,8).Value2 = 2.370349
$ws.Cells.Item(3,9).Value2 = 1
$ws.Cells.Item(3,10).Value2 = 1
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 71.80093133333332
$ws.Cells.Item(3,14).Value2 = 215.402794
$ws.Cells.Item(3,15).Value2 = 0.7804410658008428
$ws.Cells.Item(3,16).Value2 = 0.7804410658008428
$ws.Cells.Item(3,17).Value2 = 56.73108859501177
$ws.Cells.Item(3,18).Value2 = 510.5797973551059
$ws.Cells.Item(3,19).Value2 = 0.7804410658008428
$ws.Cells.Item(3,20).Value2 = 0.7804410658008428

# --- Row 4 (was old row 3, shifted down by the insert; values updated) ---
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Tff3"
$ws.Cells.Item(4,3).Value2 = "Ackr3"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 0.7901163333333333
$ws.Cells.Item(4,8).Value2 = 2.370349
$ws.Cells.Item(4,9).Value2 = 1
$ws.Cells.Item(4,10).Value2 = 1
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 10.51115833333333
$ws.Cells.Item(4,14).Value2 = 31.533475
$ws.Cells.Item(4,15).Value2 = 0.1142511588656749
$ws.Cells.Item(4,16).Value2 = 0.1142511588656749
$ws.Cells.Item(4,17).Value2 = 8.305037881419445
$ws.Cells.Item(4,18).Value2 = 74.74534093277501
$ws.Cells.Item(4,19).Value2 = 0.1142511588656749
$ws.Cells.Item(4,20).Value2 = 0.1142511588656749
